$wb = $excel.ActiveWorkbook

# --- Sheet "1963" (sheet1) ---
$ws1 = $wb.Worksheets.Item("1963")
$ws1.Cells.Item(2, 4).Value = 861127
$ws1.Cells.Item(2, 5).Value = 239614
$ws1.Cells.Item(3, 4).Value = 847994
$ws1.Cells.Item(3, 5).Value = 163769
$ws1.Cells.Item(4, 4).Value = 212415
$ws1.Cells.Item(4, 5).Value = 51183
$ws1.Cells.Item(5, 4).Value = 309215
$ws1.Cells.Item(5, 5).Value = 92444
$ws1.Cells.Item(6, 4).Value = 594537
$ws1.Cells.Item(6, 5).Value = 143897
$ws1.Cells.Item(7, 4).Value = 107573
$ws1.Cells.Item(7, 5).Value = 32425
$ws1.Cells.Item(8, 4).Value = 50836
$ws1.Cells.Item(9, 4).Value = 135126
$ws1.Cells.Item(9, 5).Value = 48483
$ws1.Cells.Item(10, 4).Value = 178337
$ws1.Cells.Item(10, 5).Value = 53445
$ws1.Cells.Item(11, 4).Value = 189607
$ws1.Cells.Item(11, 5).Value = 40620
$ws1.Cells.Item(12, 4).Value = 47898
$ws1.Cells.Item(12, 5).Value = 16036
$ws1.Cells.Item(13, 4).Value = 135658
$ws1.Cells.Item(13, 5).Value = 31042
$ws1.Cells.Item(14, 4).Value = 110567
$ws1.Cells.Item(14, 5).Value = 26334
$ws1.Cells.Item(15, 4).Value = 333367
$ws1.Cells.Item(15, 5).Value = 104072
$ws1.Cells.Item(16, 4).Value = 30464
$ws1.Cells.Item(16, 5).Value = 14349
$ws1.Cells.Item(18, 4).Value = 22665
$ws1.Cells.Item(19, 4).Value = 56897
$ws1.Cells.Item(19, 5).Value = 17258
$ws1.Cells.Item(20, 4).Value = 14661
$ws1.Cells.Item(21, 4).Value = 148876
$ws1.Cells.Item(21, 5).Value = 34255
$ws1.Cells.Item(22, 4).Value = 197548
$ws1.Cells.Item(22, 5).Value = 58691
$ws1.Cells.Item(23, 4).Value = 221475
$ws1.Cells.Item(23, 5).Value = 82157
$ws1.Cells.Item(24, 4).Value = 571095
$ws1.Cells.Item(24, 5).Value = 145495
$ws1.Cells.Item(25, 4).Value = 10492
$ws1.Cells.Item(25, 5).Value = 2701
$ws1.Cells.Item(26, 4).Value = 17266
$ws1.Cells.Item(27, 4).Value = 61754
$ws1.Cells.Item(28, 4).Value = 96526
$ws1.Cells.Item(28, 5).Value = 26462
$ws1.Cells.Item(29, 4).Value = 157559
$ws1.Cells.Item(29, 5).Value = 51779
$ws1.Cells.Item(30, 4).Value = 198007
$ws1.Cells.Item(30, 5).Value = 44521
$ws1.Cells.Item(31, 4).Value = 128360
$ws1.Cells.Item(31, 5).Value = 37853
$ws1.Cells.Item(32, 4).Value = 33190
$ws1.Cells.Item(33, 4).Value = 43758
$ws1.Cells.Item(33, 5).Value = 18125
$ws1.Cells.Item(34, 4).Value = 32577
$ws1.Cells.Item(34, 5).Value = 9678
$ws1.Cells.Item(35, 4).Value = 180070
$ws1.Cells.Item(35, 5).Value = 43680
$ws1.Cells.Item(36, 4).Value = 2176767
$ws1.Cells.Item(36, 5).Value = 443722
$ws1.Cells.Item(37, 4).Value = 17045
$ws1.Cells.Item(37, 5).Value = 6680
$ws1.Cells.Item(38, 4).Value = 44194
$ws1.Cells.Item(38, 5).Value = 14204
$ws1.Cells.Item(39, 4).Value = 46061
$ws1.Cells.Item(39, 5).Value = 16971
$ws1.Cells.Item(40, 4).Value = 7616
$ws1.Cells.Item(40, 5).Value = 2322
$ws1.Cells.Item(43, 4).Value = 95355
$ws1.Cells.Item(43, 5).Value = 31566
$ws1.Cells.Item(45, 4).Value = 35172
$ws1.Cells.Item(45, 5).Value = 12473
$ws1.Cells.Item(46, 4).Value = 239735
$ws1.Cells.Item(46, 5).Value = 84551
$ws1.Cells.Item(47, 4).Value = 245556
$ws1.Cells.Item(47, 5).Value = 71542
$ws1.Cells.Item(48, 4).Value = 57402
$ws1.Cells.Item(48, 5).Value = 22233
$ws1.Cells.Item(49, 4).Value = 23553
$ws1.Cells.Item(49, 5).Value = 8932
$ws1.Cells.Item(50, 4).Value = 925191
$ws1.Cells.Item(50, 5).Value = 208741
$ws1.Cells.Item(52, 4).Value = 49928
$ws1.Cells.Item(52, 5).Value = 14076
$ws1.Cells.Item(53, 4).Value = 30669
$ws1.Cells.Item(53, 5).Value = 8317
$ws1.Cells.Item(54, 4).Value = 147637
$ws1.Cells.Item(54, 5).Value = 53066
$ws1.Cells.Item(55, 4).Value = 17164
$ws1.Cells.Item(56, 4).Value = 168173
$ws1.Cells.Item(56, 5).Value = 76282
$ws1.Cells.Item(57, 4).Value = 249770
$ws1.Cells.Item(57, 5).Value = 89611
$ws1.Cells.Item(58, 4).Value = 51587
$ws1.Cells.Item(58, 5).Value = 18136
$ws1.Cells.Item(59, 4).Value = 759047
$ws1.Cells.Item(59, 5).Value = 181260
$ws1.Cells.Item(60, 4).Value = 25773
$ws1.Cells.Item(60, 5).Value = 8706
$ws1.Cells.Item(61, 4).Value = 36894
$ws1.Cells.Item(61, 5).Value = 11647
$ws1.Cells.Item(62, 4).Value = 129759
$ws1.Cells.Item(62, 5).Value = 44626
$ws1.Cells.Item(63, 4).Value = 37140
$ws1.Cells.Item(63, 5).Value = 12328
$ws1.Cells.Item(64, 4).Value = 41882
$ws1.Cells.Item(64, 5).Value = 13219
$ws1.Cells.Item(65, 4).Value = 323086
$ws1.Cells.Item(65, 5).Value = 103555
$ws1.Cells.Item(67, 4).Value = 33032
$ws1.Cells.Item(68, 4).Value = 144710
$ws1.Cells.Item(68, 5).Value = 52918
$ws1.Cells.Item(69, 4).Value = 14532
$ws1.Cells.Item(69, 5).Value = 4618
$ws1.Cells.Item(70, 4).Value = 3977561
$ws1.Cells.Item(70, 5).Value = 783284
$ws1.Cells.Item(71, 4).Value = 170897
$ws1.Cells.Item(71, 5).Value = 46843
$ws1.Cells.Item(73, 4).Value = 69925
$ws1.Cells.Item(73, 5).Value = 22807
$ws1.Cells.Item(74, 4).Value = 10074
$ws1.Cells.Item(74, 5).Value = 2986
$ws1.Cells.Item(76, 4).Value = 23254
$ws1.Cells.Item(76, 5).Value = 6372
$ws1.Cells.Item(77, 4).Value = 42628
$ws1.Cells.Item(77, 5).Value = 15558
$ws1.Cells.Item(78, 4).Value = 25946
$ws1.Cells.Item(79, 4).Value = 348956
$ws1.Cells.Item(79, 5).Value = 95276
$ws1.Cells.Item(80, 4).Value = 9329
$ws1.Cells.Item(80, 5).Value = 5653
$ws1.Cells.Item(81, 4).Value = 15431
$ws1.Cells.Item(81, 5).Value = 5134
$ws1.Cells.Item(82, 4).Value = 87572
$ws1.Cells.Item(82, 5).Value = 32790
$ws1.Cells.Item(83, 4).Value = 3479
$ws1.Cells.Item(85, 4).Value = 9202
$ws1.Cells.Item(85, 5).Value = 4059
$ws1.Cells.Item(86, 4).Value = 112527
$ws1.Cells.Item(86, 5).Value = 45329
$ws1.Cells.Item(87, 4).Value = 304743
$ws1.Cells.Item(87, 5).Value = 94996
$ws1.Cells.Item(88, 4).Value = 467329
$ws1.Cells.Item(88, 5).Value = 145669
$ws1.Cells.Item(89, 4).Value = 90309
$ws1.Cells.Item(89, 5).Value = 33010
$ws1.Cells.Item(93, 4).Value = 9638
$ws1.Cells.Item(95, 4).Value = 39247
$ws1.Cells.Item(95, 5).Value = 13521
$ws1.Cells.Item(96, 4).Value = 6339
$ws1.Cells.Item(96, 5).Value = 2168
$ws1.Cells.Item(97, 4).Value = 52038
$ws1.Cells.Item(97, 5).Value = 14428
$ws1.Cells.Item(98, 4).Value = 77842
$ws1.Cells.Item(98, 5).Value = 28247
$ws1.Cells.Item(99, 4).Value = 5800
$ws1.Cells.Item(99, 5).Value = 2917
$ws1.Cells.Item(100, 4).Value = 14725
$ws1.Cells.Item(100, 5).Value = 4467
$ws1.Cells.Item(101, 4).Value = 12006
$ws1.Cells.Item(101, 5).Value = 4174

# --- Sheet "1964" (sheet2) ---
$ws2 = $wb.Worksheets.Item("1964")
$ws2.Cells.Item(2, 4).Value = 4109365
$ws2.Cells.Item(2, 5).Value = 696850
$ws2.Cells.Item(3, 4).Value = 3401013
$ws2.Cells.Item(3, 5).Value = 599010
$ws2.Cells.Item(4, 4).Value = 160288
$ws2.Cells.Item(4, 5).Value = 46237
$ws2.Cells.Item(5, 4).Value = 1459708
$ws2.Cells.Item(5, 5).Value = 364031
$ws2.Cells.Item(6, 4).Value = 1805782
$ws2.Cells.Item(6, 5).Value = 400416
$ws2.Cells.Item(7, 4).Value = 499767
$ws2.Cells.Item(7, 5).Value = 154442
$ws2.Cells.Item(8, 4).Value = 773267
$ws2.Cells.Item(8, 5).Value = 247613
$ws2.Cells.Item(9, 4).Value = 27692
$ws2.Cells.Item(9, 5).Value = 8196
$ws2.Cells.Item(10, 4).Value = 67721
$ws2.Cells.Item(10, 5).Value = 14003
$ws2.Cells.Item(11, 4).Value = 741198
$ws2.Cells.Item(11, 5).Value = 206424
$ws2.Cells.Item(12, 4).Value = 130516
$ws2.Cells.Item(12, 5).Value = 39089
$ws2.Cells.Item(13, 4).Value = 50135
$ws2.Cells.Item(13, 5).Value = 16116
$ws2.Cells.Item(14, 4).Value = 3939016
$ws2.Cells.Item(14, 5).Value = 672797
$ws2.Cells.Item(15, 4).Value = 4026316
$ws2.Cells.Item(15, 5).Value = 703418
$ws2.Cells.Item(16, 4).Value = 507747
$ws2.Cells.Item(16, 5).Value = 152167
$ws2.Cells.Item(17, 4).Value = 1823523
$ws2.Cells.Item(17, 5).Value = 357287
$ws2.Cells.Item(18, 4).Value = 84146
$ws2.Cells.Item(18, 5).Value = 28395
$ws2.Cells.Item(21, 4).Value = 677670
$ws2.Cells.Item(21, 5).Value = 191105
$ws2.Cells.Item(22, 4).Value = 353842
$ws2.Cells.Item(22, 5).Value = 112147
$ws2.Cells.Item(23, 4).Value = 46197
$ws2.Cells.Item(23, 5).Value = 13696
$ws2.Cells.Item(24, 4).Value = 90083
$ws2.Cells.Item(24, 5).Value = 21122
$ws2.Cells.Item(25, 4).Value = 31652
$ws2.Cells.Item(25, 5).Value = 10446
$ws2.Cells.Item(26, 4).Value = 9484
$ws2.Cells.Item(26, 5).Value = 2528
$ws2.Cells.Item(27, 4).Value = 52679
$ws2.Cells.Item(27, 5).Value = 16133
$ws2.Cells.Item(28, 4).Value = 40138
$ws2.Cells.Item(28, 5).Value = 12564
$ws2.Cells.Item(29, 4).Value = 76949
$ws2.Cells.Item(29, 5).Value = 31755
$ws2.Cells.Item(30, 4).Value = 409854
$ws2.Cells.Item(30, 5).Value = 96353
$ws2.Cells.Item(31, 4).Value = 176762
$ws2.Cells.Item(31, 5).Value = 55243
$ws2.Cells.Item(32, 4).Value = 122288
$ws2.Cells.Item(32, 5).Value = 33528
$ws2.Cells.Item(33, 4).Value = 7175
$ws2.Cells.Item(34, 4).Value = 989074
$ws2.Cells.Item(34, 5).Value = 282739
$ws2.Cells.Item(36, 4).Value = 185355
$ws2.Cells.Item(36, 5).Value = 49012
$ws2.Cells.Item(37, 4).Value = 431052
$ws2.Cells.Item(37, 5).Value = 94139
$ws2.Cells.Item(38, 4).Value = 740142
$ws2.Cells.Item(38, 5).Value = 209205
$ws2.Cells.Item(39, 4).Value = 5757895
$ws2.Cells.Item(39, 5).Value = 1003651
$ws2.Cells.Item(40, 4).Value = 20720
$ws2.Cells.Item(40, 5).Value = 6999
$ws2.Cells.Item(41, 4).Value = 3227729
$ws2.Cells.Item(41, 5).Value = 498851
$ws2.Cells.Item(42, 4).Value = 57971
$ws2.Cells.Item(42, 5).Value = 23651
$ws2.Cells.Item(43, 4).Value = 5848
$ws2.Cells.Item(43, 5).Value = 2308
$ws2.Cells.Item(44, 4).Value = 21605
$ws2.Cells.Item(44, 5).Value = 7505
$ws2.Cells.Item(45, 4).Value = 123540
$ws2.Cells.Item(45, 5).Value = 40140
$ws2.Cells.Item(46, 4).Value = 38596
$ws2.Cells.Item(46, 5).Value = 10312
$ws2.Cells.Item(47, 4).Value = 113681
$ws2.Cells.Item(47, 5).Value = 27599
$ws2.Cells.Item(48, 4).Value = 80565
$ws2.Cells.Item(49, 4).Value = 24988
$ws2.Cells.Item(49, 5).Value = 10345
$ws2.Cells.Item(50, 4).Value = 108281
$ws2.Cells.Item(50, 5).Value = 36105
$ws2.Cells.Item(51, 4).Value = 325073
$ws2.Cells.Item(51, 5).Value = 88354
$ws2.Cells.Item(53, 4).Value = 4106123
$ws2.Cells.Item(53, 5).Value = 694667
$ws2.Cells.Item(54, 4).Value = 111578
$ws2.Cells.Item(54, 5).Value = 32590
$ws2.Cells.Item(55, 4).Value = 14764
$ws2.Cells.Item(55, 5).Value = 6222
$ws2.Cells.Item(56, 4).Value = 1285093
$ws2.Cells.Item(56, 5).Value = 245258
$ws2.Cells.Item(57, 4).Value = 71190
$ws2.Cells.Item(57, 5).Value = 32376
$ws2.Cells.Item(58, 4).Value = 540944
$ws2.Cells.Item(58, 5).Value = 172787
$ws2.Cells.Item(59, 4).Value = 17350
$ws2.Cells.Item(59, 5).Value = 4609
$ws2.Cells.Item(60, 4).Value = 28520
$ws2.Cells.Item(60, 5).Value = 11989
$ws2.Cells.Item(61, 4).Value = 297038
$ws2.Cells.Item(61, 5).Value = 97086
$ws2.Cells.Item(62, 4).Value = 16987
$ws2.Cells.Item(62, 5).Value = 5996
$ws2.Cells.Item(63, 4).Value = 67242
$ws2.Cells.Item(63, 5).Value = 24083
$ws2.Cells.Item(64, 4).Value = 11909
$ws2.Cells.Item(64, 5).Value = 3531
$ws2.Cells.Item(65, 4).Value = 57681
$ws2.Cells.Item(68, 4).Value = 9799
$ws2.Cells.Item(68, 5).Value = 2847
$ws2.Cells.Item(69, 4).Value = 20033
$ws2.Cells.Item(69, 5).Value = 6089
$ws2.Cells.Item(70, 4).Value = 539865
$ws2.Cells.Item(70, 5).Value = 160249
$ws2.Cells.Item(72, 4).Value = 553458
$ws2.Cells.Item(72, 5).Value = 179622
$ws2.Cells.Item(73, 4).Value = 174980
$ws2.Cells.Item(73, 5).Value = 55502
$ws2.Cells.Item(75, 4).Value = 244683
$ws2.Cells.Item(75, 5).Value = 80379
$ws2.Cells.Item(77, 4).Value = 1001243
$ws2.Cells.Item(77, 5).Value = 211458
$ws2.Cells.Item(78, 4).Value = 14335
$ws2.Cells.Item(78, 5).Value = 5665
$ws2.Cells.Item(79, 4).Value = 27411
$ws2.Cells.Item(79, 5).Value = 11454
$ws2.Cells.Item(80, 4).Value = 3112257
$ws2.Cells.Item(80, 5).Value = 698045
$ws2.Cells.Item(81, 4).Value = 12750
$ws2.Cells.Item(82, 4).Value = 58588
$ws2.Cells.Item(82, 5).Value = 27829
$ws2.Cells.Item(83, 4).Value = 18145
$ws2.Cells.Item(84, 4).Value = 85866
$ws2.Cells.Item(84, 5).Value = 29173
$ws2.Cells.Item(85, 4).Value = 10381
$ws2.Cells.Item(86, 4).Value = 37475
$ws2.Cells.Item(86, 5).Value = 12898
$ws2.Cells.Item(87, 4).Value = 23964
$ws2.Cells.Item(87, 5).Value = 10079
$ws2.Cells.Item(88, 4).Value = 94596
$ws2.Cells.Item(88, 5).Value = 36000
$ws2.Cells.Item(90, 4).Value = 74050
$ws2.Cells.Item(90, 5).Value = 28653
$ws2.Cells.Item(92, 4).Value = 31946
$ws2.Cells.Item(92, 5).Value = 9238
$ws2.Cells.Item(93, 4).Value = 37296
$ws2.Cells.Item(93, 5).Value = 11123
$ws2.Cells.Item(94, 4).Value = 13223
$ws2.Cells.Item(95, 4).Value = 29971
$ws2.Cells.Item(95, 5).Value = 10566
$ws2.Cells.Item(96, 4).Value = 2527641
$ws2.Cells.Item(96, 5).Value = 420321
$ws2.Cells.Item(99, 4).Value = 14075
$ws2.Cells.Item(99, 5).Value = 4954
$ws2.Cells.Item(100, 4).Value = 1030621
$ws2.Cells.Item(100, 5).Value = 287586
$ws2.Cells.Item(101, 4).Value = 296391
$ws2.Cells.Item(101, 5).Value = 87838
